# CERN_QTR_FIN.xlsx quarterly financials update
#
# The source workbook carries 8 quarters of data in columns D:K (most
# recent quarter first, in column D). This commit adds the two most
# recent quarters (Q4 2018 / Q3 2018) as new columns D:E, pushing all of
# the existing quarters two columns to the right (old D:K -> new F:M).
#
# Row 91 ("Capital Expenditures") also carries restated historical
# figures for the five quarters that land in F:J after the shift, so
# those five cells get explicit corrected values as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns before D; everything in D:K moves to F:M.
$ws.Columns("D:E").Insert()

# 2) The new D:E columns come in with default/general formatting; copy
#    the number formats / fonts / alignment from F:G (which is what used
#    to be D:E before the insert) so D:E look like the rest of the table.
$ws.Range("F5:G102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 3) Populate the two new quarter columns (D = newest quarter, E = next).
#    Keyed by row number -> @(D-value, E-value).
$newQuarterData = @{
    7   = @(43463, 43372)
    8   = @(1365700, 1340100)
    9   = @(237000, 230300)
    10  = @(1128700, 1109800)
    12  = @(126800, 118900)
    13  = @(0, 0)
    14  = @(0, 0)
    15  = @(76100, 75000)
    17  = @(1201500, 1132900)
    18  = @(164200, 207200)
    20  = @(7600, 6900)
    21  = @(340600, 375200)
    22  = @(0, 0)
    23  = @(171800, 214100)
    24  = @(40500, 44700)
    25  = @(0, 0)
    26  = @(131300, 169400)
    27  = @(131300, 169400)
    28  = @(0, 0)
    29  = @("NA", "NA")
    30  = @(0, 0)
    31  = @(0, 0)
    32  = @(-7600, -6900)
    33  = @(131300, 169400)
    34  = @(0, 0)
    35  = @(131300, 169400)
    38  = @(43463, 43372)
    41  = @(374100, 498600)
    42  = @(401300, 314900)
    43  = @(1183500, 1210600)
    44  = @(25000, 24700)
    45  = @(334900, 337900)
    46  = @(2318800, 2386800)
    47  = @(300000, 338200)
    48  = @(1743600, 1697200)
    49  = @(2147400, 2151000)
    50  = @(0, 0)
    51  = @(0, 0)
    52  = @(198900, 212100)
    53  = @(0, 0)
    54  = @(6708600, 6785400)
    57  = @(293500, 258400)
    58  = @(4900, 2300)
    59  = @(664200, 646300)
    60  = @(962700, 907000)
    61  = @(438800, 438800)
    62  = @(378800, 373900)
    63  = @(0, 0)
    64  = @(0, 0)
    65  = @(0, 0)
    66  = @(1780200, 1719700)
    68  = @(0, 0)
    69  = @(0, 0)
    70  = @(0, 0)
    71  = @(0, 0)
    72  = @(5576500, 5445200)
    73  = @(0, 0)
    74  = @(0, 0)
    75  = @(0, 0)
    76  = @(4928400, 5065700)
    77  = @(0, 0)
    80  = @(43463, 43372)
    81  = @(131300, 169400)
    83  = @(168800, 161100)
    84  = @(0, 0)
    85  = @(0, 0)
    86  = @(0, 0)
    87  = @(0, 0)
    88  = @(0, 0)
    89  = @(406900, 338500)
    91  = @(-141000, -117000)
    92  = @(0, 0)
    93  = @(0, 0)
    94  = @(-266800, -350900)
    96  = @(0, 0)
    97  = @(0, 0)
    98  = @(0, 0)
    99  = @(0, 0)
    100 = @(-264100, 4600)
    101 = @(-500, -4500)
    102 = @(-124500, -12400)
}

foreach ($row in $newQuarterData.Keys) {
    $pair = $newQuarterData[$row]
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# 4) Row 91 ("Capital Expenditures") also has restated values for the
#    five quarters that now sit in F:J (previously D:H) - overwrite them.
$row91Restated = @{
    6 = -109300   # F91
    7 = -79700    # G91
    8 = -99700    # H91
    9 = -73000    # I91
    10 = -101300  # J91
}
foreach ($col in $row91Restated.Keys) {
    $ws.Cells.Item(91, $col).Value = $row91Restated[$col]
}
